$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''21.641.80'
$ws.Range("E2").Value = '''  -1.98%  '
$ws.Range("D3").Value = '''1.533.83'
$ws.Range("E3").Value = '''  -1.53%  '
$ws.Range("D5").Value = '''1.001'
$ws.Range("E5").Value = '''  +0.11%  '
$ws.Range("D6").Value = '''288.33'
$ws.Range("E6").Value = '''  +0.21%  '
$ws.Range("D7").Value = '''0.3935'
$ws.Range("E7").Value = '''  +1.79%  '
$ws.Range("D8").Value = '''0.3166'
$ws.Range("E8").Value = '''  -2.42%  '
$ws.Range("D9").Value = '''42.45'
$ws.Range("E9").Value = '''  -1.54%  '
$ws.Range("D10").Value = '''0.07166'
$ws.Range("D11").Value = '''1.049'
$ws.Range("E11").Value = '''  -6.69%  '
$ws.Range("E12").Value = '''  +0.11%  '
$ws.Range("D13").Value = '''5.665'
$ws.Range("E13").Value = '''  -0.69%  '
$ws.Range("D14").Value = '''18.55'
$ws.Range("E14").Value = '''  -4.27%  '
$ws.Range("D15").Value = '''6.581'
$ws.Range("E15").Value = '''  -3.37%  '
$ws.Range("D16").Value = '''1.531.46'
$ws.Range("E16").Value = '''  -1.52%  '
$ws.Range("D17").Value = '''0.00001089'
$ws.Range("E17").Value = '''  -2.78%  '
$ws.Range("D18").Value = '''0.06594'
$ws.Range("E18").Value = '''  -0.23%  '
$ws.Range("D19").Value = '''83.69'
$ws.Range("E19").Value = '''  -1.88%  '
$ws.Range("D20").Value = '''1.001'
$ws.Range("E20").Value = '''  +0.17%  '
$ws.Range("D21").Value = '''6.113'
$ws.Range("E21").Value = '''  -4.52%  '
$ws.Range("D22").Value = '''15.44'
$ws.Range("E22").Value = '''  -3.48%  '
$ws.Range("D23").Value = '''10.74'
$ws.Range("E23").Value = '''  -6.46%  '
$ws.Range("D24").Value = '''2.346'
$ws.Range("E24").Value = '''  +0.59%  '
$ws.Range("D25").Value = '''21.642.33'
$ws.Range("E25").Value = '''  -1.98%  '
$ws.Range("D26").Value = '''2.348'
$ws.Range("E26").Value = '''  -8.12%  '
$ws.Range("D27").Value = '''149.26'
$ws.Range("E27").Value = '''  -0.24%  '
$ws.Range("D28").Value = '''18.30'
$ws.Range("E28").Value = '''  -3.06%  '
$ws.Range("D29").Value = '''4.845'
$ws.Range("E29").Value = '''  -0.47%  '
$ws.Range("D30").Value = '''1.746.63'
$ws.Range("E30").Value = '''  +0.91%  '
$ws.Range("D31").Value = '''117.08'
$ws.Range("E31").Value = '''  -3.23%  '
$ws.Range("D32").Value = '''6.005'
$ws.Range("E32").Value = '''  +2.36%  '
$ws.Range("D33").Value = '''0.9398'
$ws.Range("E33").Value = '''  -15.66%  '
$ws.Range("D34").Value = '''0.08143'
$ws.Range("E34").Value = '''  -0.62%  '
$ws.Range("D35").Value = '''8.499'
$ws.Range("E35").Value = '''  -8.46%  '
$ws.Range("D36").Value = '''5.147'
$ws.Range("E36").Value = '''  -1.60%  '
$ws.Range("D37").Value = '''0.06013'
$ws.Range("E37").Value = '''  -3.68%  '
$ws.Range("D38").Value = '''0.02215'
$ws.Range("E38").Value = '''  -3.80%  '
$ws.Range("D39").Value = '''1.456'
$ws.Range("E39").Value = '''  -14.34%  '
$ws.Range("E40").Value = '''  -4.26%  '
$ws.Range("D41").Value = '''1.178'
$ws.Range("E41").Value = '''  -3.50%  '
$ws.Range("D42").Value = '''10.93'
$ws.Range("E42").Value = '''  +0.39%  '
$ws.Range("D43").Value = '''0.9997'
$ws.Range("D44").Value = '''0.5750'
$ws.Range("E44").Value = '''  -3.45%  '
$ws.Range("D45").Value = '''12.96'
$ws.Range("E45").Value = '''  -4.32%  '
$ws.Range("D46").Value = '''3.707'
$ws.Range("E46").Value = '''  -0.29%  '
$ws.Range("D47").Value = '''0.5492'
$ws.Range("E47").Value = '''  -4.61%  '
$ws.Range("D48").Value = '''1.165'
$ws.Range("E48").Value = '''  +0.41%  '
$ws.Range("E49").Value = '''  -2.81%  '
$ws.Range("D50").Value = '''116.04'
$ws.Range("E50").Value = '''  -2.69%  '
$ws.Range("D51").Value = '''0.06689'
$ws.Range("E51").Value = '''  -2.93%  '
